$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 2022

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = 92.960099223795225

$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = 96.03949422949897

$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = 91.012153547624152

$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)

$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Value = 94.391087218067838

$ws.Range("G11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = 91.76755842559642

$ws.Range("G12").Copy()
$ws.Range("H12").PasteSpecial(-4122)

$ws.Range("G13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = 92.942689638142156

$ws.Range("G14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = 86.897877953385489

$ws.Range("G15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 96.500794494289821

$ws.Range("G16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = 94.135975315309977

$ws.Range("G17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H17").Value = 89.456106196597958

$ws.Range("G18").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H18").Value = 94.270923428904894

$ws.Range("G19").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value = 97.027480110114013

$ws.Range("G20").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = 98.077227596867303

$ws.Range("G21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H21").Value = 90.983384827072243

$ws.Range("G22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("G23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").Value = 90.468970496790078

$ws.Range("G24").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("H24").Value = 95.809965597614095

$ws.Range("G25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value = 88.221110530662017

$ws.Range("G26").Copy()
$ws.Range("H26").PasteSpecial(-4122)

$ws.Range("G27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = 69.811292606515579

$ws.Range("G28").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = 85.757158930558518

$ws.Range("G29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 93.032103866435918

$ws.Range("G30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = 97.325262246493097

$ws.Range("G31").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = 98.908492141713779

$ws.Range("G32").Copy()
$ws.Range("H32").PasteSpecial(-4122)

$ws.Range("G33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = 91.968006037496949

$ws.Range("G34").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("H34").Value = 91.809335747904541

$ws.Range("G35").Copy()
$ws.Range("H35").PasteSpecial(-4122)
$ws.Range("H35").Value = 91.27524653351901

$ws.Range("G36").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("H36").Value = 94.397288657466234

$ws.Range("G37").Copy()
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("H37").Value = 96.740699993405215

$excel.CutCopyMode = $false
$ws.Range("I4").Select()
